$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 11000.214
$ws.Range("I74").Value = 11000.214
$ws.Range("K74").Value = 11000.214
$ws.Range("M74").Value = -10064.214
$ws.Range("H77").Value = 11000.214
$ws.Range("I77").Value = 11000.214
$ws.Range("K77").Value = 55001.07
$ws.Range("M77").Value = -50321.07
$ws.Range("H103").Value = 853.3570999999999
$ws.Range("I103").Value = 493
$ws.Range("J103").Value = 951.63635
$ws.Range("K103").Value = 1479
$ws.Range("L103").Value = 2854.90905
$ws.Range("M103").Value = -893
$ws.Range("N103").Value = -4026.90905
$ws.Range("H135").Value = 7861.643
$ws.Range("I135").Value = 5587.8335
$ws.Range("J135").Value = 9567
$ws.Range("K135").Value = 50290.5015
$ws.Range("L135").Value = 86103
$ws.Range("M135").Value = -47755.5015
$ws.Range("N135").Value = -91173

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2672.5
$ws.Range("I63").Value = 2672.5
$ws.Range("K63").Value = 2672.5
$ws.Range("M63").Value = -1986.5
$ws.Range("H66").Value = 2672.5
$ws.Range("I66").Value = 2672.5
$ws.Range("K66").Value = 13362.5
$ws.Range("M66").Value = -9930.5
$ws.Range("H74").Value = 4998.4917
$ws.Range("I74").Value = 2317.9387
$ws.Range("K74").Value = 2317.9387
$ws.Range("M74").Value = -1443.9387
$ws.Range("H77").Value = 4998.4917
$ws.Range("I77").Value = 2317.9387
$ws.Range("K77").Value = 11589.6935
$ws.Range("M77").Value = -7221.693500000001
$ws.Range("H88").Value = 2467.1333
$ws.Range("J88").Value = 3751.75
$ws.Range("L88").Value = 3751.75
$ws.Range("N88").Value = -4563.75
$ws.Range("H91").Value = 2467.1333
$ws.Range("J91").Value = 3751.75
$ws.Range("L91").Value = 3751.75
$ws.Range("N91").Value = -6559.75
$ws.Range("H134").Value = 200000
$ws.Range("J134").Value = 200000
$ws.Range("L134").Value = 200000
$ws.Range("N134").Value = -210140

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6248.5
$ws.Range("I86").Value = 6664.6665
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 6664.6665
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -5541.6665
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 6248.5
$ws.Range("I89").Value = 6664.6665
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 33323.3325
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -27707.3325
$ws.Range("N89").Value = -36232
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 5163
$ws.Range("I134").Value = 5163
$ws.Range("K134").Value = 15489
$ws.Range("M134").Value = -12954

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 38500
$ws.Range("J70").Value = 38500
$ws.Range("L70").Value = 38500
$ws.Range("N70").Value = -39130
$ws.Range("H73").Value = 38500
$ws.Range("J73").Value = 38500
$ws.Range("L73").Value = 38500
$ws.Range("N73").Value = -40684
$ws.Range("H94").Value = 6678.722
$ws.Range("I94").Value = 17122.334
$ws.Range("K94").Value = 17122.334
$ws.Range("M94").Value = -16671.334
$ws.Range("H99").Value = 18949.834
$ws.Range("I99").Value = 22048.4
$ws.Range("J99").Value = 3457
$ws.Range("K99").Value = 22048.4
$ws.Range("L99").Value = 3457
$ws.Range("M99").Value = -20550.4
$ws.Range("N99").Value = -6453
$ws.Range("H126").Value = 18949.834
$ws.Range("I126").Value = 22048.4
$ws.Range("J126").Value = 3457
$ws.Range("K126").Value = 66145.20000000001
$ws.Range("L126").Value = 10371
$ws.Range("M126").Value = -63675.20000000001
$ws.Range("N126").Value = -15311
$ws.Range("H134").Value = 2252.7646
$ws.Range("I134").Value = 2252.7646
$ws.Range("K134").Value = 6758.293799999999
$ws.Range("M134").Value = -4223.293799999999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5336066
$ws.Range("I4").Value = 8889888
$ws.Range("K4").Value = 26669664
$ws.Range("M4").Value = -26669552
$ws.Range("H7").Value = 62.8
$ws.Range("I7").Value = 40.666668
$ws.Range("J7").Value = 96
$ws.Range("K7").Value = 122.000004
$ws.Range("L7").Value = 288
$ws.Range("M7").Value = -10.000004
$ws.Range("N7").Value = -512
$ws.Range("H9").Value = 201000.05
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 201000.05
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 603000.1499999999
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -603448.1499999999
$ws.Range("H17").Value = 1682
$ws.Range("I17").Value = 1599.7142
$ws.Range("J17").Value = 1874
$ws.Range("K17").Value = 4799.142599999999
$ws.Range("L17").Value = 5622
$ws.Range("M17").Value = -4630.142599999999
$ws.Range("N17").Value = -5960
$ws.Range("H50").Value = 116588.49
$ws.Range("I50").Value = 309.41177
$ws.Range("J50").Value = 192617.11
$ws.Range("K50").Value = 928.23531
$ws.Range("L50").Value = 577851.33
$ws.Range("M50").Value = -447.23531
$ws.Range("N50").Value = -578813.33
$ws.Range("H53").Value = 116588.49
$ws.Range("I53").Value = 309.41177
$ws.Range("J53").Value = 192617.11
$ws.Range("K53").Value = 928.23531
$ws.Range("L53").Value = 577851.33
$ws.Range("M53").Value = -447.23531
$ws.Range("N53").Value = -578813.33
$ws.Range("H59").Value = 1500
$ws.Range("I59").Value = 1500
$ws.Range("K59").Value = 4500
$ws.Range("M59").Value = -3960

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 6997.3
$ws.Range("I21").Value = 6108.6665
$ws.Range("K21").Value = 6108.6665
$ws.Range("M21").Value = -5935.6665
$ws.Range("H24").Value = 1013333.3
$ws.Range("J24").Value = 1013333.3
$ws.Range("L24").Value = 1013333.3
$ws.Range("N24").Value = -1013679.3
$ws.Range("H30").Value = 6997.3
$ws.Range("I30").Value = 6108.6665
$ws.Range("K30").Value = 6108.6665
$ws.Range("M30").Value = -6003.6665
$ws.Range("H69").Value = 35500
$ws.Range("J69").Value = 36000
$ws.Range("L69").Value = 36000
$ws.Range("N69").Value = -37498
$ws.Range("H72").Value = 35500
$ws.Range("J72").Value = 36000
$ws.Range("L72").Value = 108000
$ws.Range("N72").Value = -115488

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 43533.816
$ws.Range("I74").Value = 34270.25
$ws.Range("J74").Value = 48827.285
$ws.Range("K74").Value = 34270.25
$ws.Range("L74").Value = 48827.285
$ws.Range("M74").Value = -33272.25
$ws.Range("N74").Value = -50823.285
$ws.Range("H77").Value = 43533.816
$ws.Range("I77").Value = 34270.25
$ws.Range("J77").Value = 48827.285
$ws.Range("K77").Value = 102810.75
$ws.Range("L77").Value = 146481.855
$ws.Range("M77").Value = -97818.75
$ws.Range("N77").Value = -156465.855
$ws.Range("H122").Value = 4088.3809
$ws.Range("I122").Value = 2901.0625
$ws.Range("K122").Value = 8703.1875
$ws.Range("M122").Value = -6253.1875
$ws.Range("H133").Value = 82616.664
$ws.Range("J133").Value = 88777
$ws.Range("L133").Value = 88777
$ws.Range("N133").Value = -93837

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 369483.72
$ws.Range("I4").Value = 405432.2
$ws.Range("K4").Value = 405432.2
$ws.Range("M4").Value = -405319.2
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H128").Value = 59999
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H132").Value = 3348.5715
$ws.Range("I132").Value = 2638.7368
$ws.Range("K132").Value = 7916.2104
$ws.Range("M132").Value = -5386.2104
